$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "ActiveSheet name:" $ws.Name
Write-Host "Worksheets count:" $wb.Worksheets.Count
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Host "Sheet $i :" $wb.Worksheets.Item($i).Name
}
